$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-50: Wins=84, Losses=78, Ties=0
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 84   # AD
    $ws.Cells.Item($r, 31).Value = 78   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
